$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Chin Pok%CHAN%chinpok0301@gmail.com%0, Ngai Sze%WONG%candy_wong@cuhk.edu.hk%2, Ngai Sze%WONG%candy_wong@cuhk.edu.hk%0, Chi Chiu%LEUNG%ccleungpnc@netvigator.com%2, Chi Chiu%LEUNG%ccleungpnc@netvigator.com%0, Shui Shan%LEE%sslee@cuhk.edu.hk%1]"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "Oxford University Press"
$ws.Range("C3").Value = "Unknown Title"
$ws.Range("E3").Value = "[]"
$ws.Range("F3").Value = "not found"
$ws.Range("G3").Value = "N/A"
$ws.Range("I3").Value = ""
$ws.Range("E4").Value = "[Khee-Siang%Chan%NULL%0, Fu-Wen%Liang%NULL%1, Hung-Jen%Tang%NULL%1, Han Siong%Toh%NULL%1, Wen-Liang%Yu%NULL%1]"
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = "Elsevier España, S.L.U."
$ws.Range("E5").Value = "[Young June%Choe%NULL%0, Jong-Koo%Lee%NULL%2, Jong-Koo%Lee%NULL%0]"
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = "Korea Centers for Disease Control and Prevention"
$ws.Range("E6").Value = "[Ermengol%Coma Redon%NULL%0, Nuria%Mora%NULL%1, Albert%Prats-Uribe%NULL%0, Francesc%Fina Avilés%NULL%2, Francesc%Fina Avilés%NULL%0, Daniel%Prieto-Alhambra%NULL%0, Manuel%Medina%NULL%2, Manuel%Medina%NULL%0]"
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = "BMJ Publishing Group"
$ws.Range("E7").Value = "[Benjamin J%Cowling%NULL%0, Sheikh Taslim%Ali%NULL%2, Tiffany W Y%Ng%NULL%2, Tim K%Tsang%NULL%2, Julian C M%Li%NULL%2, Min Whui%Fong%NULL%2, Qiuyan%Liao%NULL%2, Mike YW%Kwan%NULL%2, So Lun%Lee%NULL%2, Susan S%Chiu%NULL%2, Joseph T%Wu%NULL%2, Peng%Wu%pengwu@hku.hk%0, Gabriel M%Leung%NULL%0]"
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = "The Author(s). Published by Elsevier Ltd."
$ws.Range("D8").Value = "
              ⬢
              We explored the possibility that public concern over COVID-19 may have impacted the seasonal trend of influenza in the northern hemisphere in winter and early spring.
"
$ws.Range("E8").Value = "[Takahiro%Itaya%NULL%0, Yuki%Furuse%NULL%1, Kazuaki%Jindai%NULL%1]"
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = "The Author(s). Published by Elsevier Ltd on behalf of International Society for Infectious Diseases."
$ws.Range("D9").Value = "At the end of 2019, an outbreak of pneumonia took place caused by a new coronavirus (SARS-CoV-2 virus), named coronavirus disease 2019 (COVID-19).
 A series of strict prevention and control measures were then implemented to reduce the spread of the epidemic.
 Influenza, another respiratory tract virus, may also respond to these measures.
 To assess the impact of these measures, we used the total number of passengers movement in mainland China from 2018 to 2020 and daily number of railway passenger flow during the 2020 Spring Festival travel rush to reflect the population movement and to analyze newly and cumulatively confirmed COVID-19 and influenza cases.
 We found that implementing the series of measures against COVID-19 mitigated both COVID-19 and influenza epidemics in China.
 Prevention and control measures for COVID-19 might be used to control respiratory tract infections to reduce the national health economic burden caused by these pathogens.
"
$ws.Range("E9").Value = "[Xiangsha%Kong%NULL%0, Feng%Liu%NULL%0, Haibo%Wang%NULL%0, Ruifeng%Yang%NULL%1, Dongbo%Chen%NULL%1, Xiaoxiao%Wang%NULL%1, Fengmin%Lu%lu.fengmin@hsc.pku.edu.cn%1, Huiying%Rao%raohuiying@pkuph.edu.cn%1, Hongsong%Chen%chenhongsong@bjmu.edu.cn%1]"
$ws.Range("I9").Value = ""
$ws.Range("J9").Value = "Elsevier"
$ws.Range("D10").Value = "Taiwan has strictly followed infection control measures to prevent spread of coronavirus disease.
 Meanwhile, nationwide surveillance data revealed drastic decreases in influenza diagnoses in outpatient departments, positivity rates of clinical specimens, and confirmed severe cases during the first 12 weeks of 2020 compared with the same period of 2019."
$ws.Range("E10").Value = "[Shu-Chen%Kuo%NULL%0, Shu-Man%Shih%NULL%1, Li-Hsin%Chien%NULL%1, Chao A.%Hsiung%NULL%1]"
$ws.Range("I10").Value = ""
$ws.Range("J10").Value = "Centers for Disease Control and Prevention"
$ws.Range("E11").Value = "[Hyunju%Lee%NULL%0, Heeyoung%Lee%NULL%1, Kyoung-Ho%Song%NULL%0, Eu Suk%Kim%NULL%0, Eu Suk%Kim%NULL%0, Jeong Su%Park%NULL%1, Jongtak%Jung%NULL%1, Soyeon%Ahn%NULL%1, Eun Kyeong%Jeong%NULL%1, Hyekyung%Park%NULL%1, Hong Bin%Kim%hbkimmd@snu.ac.kr%0]"
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = "Oxford University Press"
$ws.Range("D12").Value = "Social distancing has been adopted as one of basic protective measures against coronavirus disease 2019 (COVID-19).
 During 2019–2020 season, influenza epidemic period was exceptionally short and epidemic peak was low in comparison with previous seasons in Korea.
 Influenza epidemic pattern was bimodal in 2016–2017 and 2018–2019 seasons, however, influenza viruses have rarely been circulating in spring, 2020 in Korea.
 Although multiple factors could affect the size of influenza epidemic, extensive application of nonpharmaceutical interventions including mask wearing and social distancing in response to COVID-19 seems to be a major factor of reduced influenza epidemic.
 Social distancing measures with high feasibility and high acceptability should be implemented even if severe acute respiratory syndrome coronavirus 2 (SARS-CoV-2) vaccines are developed in the future.
 Establishment of guideline for workplace social distancing is needed and it would contribute to reduce disease burden of influenza, especially in vaccine mismatch year.
"
$ws.Range("E12").Value = "[Ji Yun%Noh%NULL%0, Hye%Seong%NULL%0, Hye%Seong%NULL%0, Jin Gu%Yoon%NULL%0, Jin Gu%Yoon%NULL%0, Joon Young%Song%NULL%0, Joon Young%Song%NULL%0, Hee Jin%Cheong%NULL%0, Hee Jin%Cheong%NULL%0, Woo Joo%Kim%NULL%0, Woo Joo%Kim%NULL%0]"
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = "The Korean Academy of Medical Sciences"
$ws.Range("D13").Value = "
              •
              Influenza virus transmission may be stopped while fighting the severe acute respiratory syndrome coronavirus 2 (SARS-CoV-2) outbreak.
"
$ws.Range("E13").Value = "[Di%Wu%NULL%0, Jianyun%Lu%NULL%1, Yanhui%Liu%NULL%1, Zhoubin%Zhang%NULL%1, Lei%Luo%NULL%1]"
$ws.Range("I13").Value = ""
$ws.Range("J13").Value = "Published by Elsevier Ltd on behalf of International Society for Infectious Diseases."
$ws.Range("E14").Value = "[Chih-Jen%Yang%NULL%0, Tun-Chieh%Chen%NULL%2, Tun-Chieh%Chen%NULL%0, Shin-Huei%Kuo%NULL%1, Min-Han%Hsieh%NULL%1, Yen-Hsu%Chen%NULL%1]"
$ws.Range("I14").Value = ""
$ws.Range("J14").Value = "Cambridge University Press"
$ws.Range("C15").Value = "Unknown Title"
$ws.Range("D15").Value = "Unknown Abstract"
$ws.Range("E15").Value = "[]"
$ws.Range("F15").Value = "not found"
$ws.Range("G15").Value = "N/A"
$ws.Range("H15").Value = "1970-01-01"
$ws.Range("I15").Value = ""
